$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.183.05'
$ws.Range("E2").Value = '  +2.89%  '
$ws.Range("D3").Value = '3.618.05'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '625.45'
$ws.Range("E5").Value = '  +2.54%  '
$ws.Range("D6").Value = '158.90'
$ws.Range("E6").Value = '  +3.29%  '
$ws.Range("D7").Value = '3.621.98'
$ws.Range("E7").Value = '  +2.46%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  +1.60%  '
$ws.Range("E10").Value = '  +3.12%  '
$ws.Range("D11").Value = '7.17'
$ws.Range("E11").Value = '  +4.73%  '
$ws.Range("D12").Value = '0.438'
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("D13").Value = '0.0000224'
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = '33.09'
$ws.Range("E14").Value = '  +3.67%  '
$ws.Range("D15").Value = '4.246.88'
$ws.Range("E15").Value = '  +2.58%  '
$ws.Range("D16").Value = '3.632.39'
$ws.Range("E16").Value = '  +2.86%  '
$ws.Range("D17").Value = '69.203.02'
$ws.Range("E17").Value = '  +3.00%  '
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").Value = '6.59'
$ws.Range("E19").Value = '  +4.10%  '
$ws.Range("D20").Value = '15.83'
$ws.Range("E20").Value = '  +2.72%  '
$ws.Range("D21").Value = '10.18'
$ws.Range("E21").Value = '  +9.88%  '
$ws.Range("D22").Value = '459.01'
$ws.Range("E22").Value = '  +2.87%  '
$ws.Range("D23").Value = '0.639'
$ws.Range("E23").Value = '  +1.23%  '
$ws.Range("D24").Value = '78.48'
$ws.Range("E24").Value = '  +0.60%  '
$ws.Range("D25").Value = '0.0000135'
$ws.Range("E25").Value = '  +10.81%  '
$ws.Range("D26").Value = '3.776.99'
$ws.Range("E26").Value = '  +2.62%  '
$ws.Range("D27").Value = '10.60'
$ws.Range("E27").Value = '  +3.51%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '9.15'
$ws.Range("E29").Value = '  +10.60%  '
$ws.Range("D30").Value = '2.61'
$ws.Range("E30").Value = '  +2.83%  '
$ws.Range("D31").Value = '1.70'
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("D33").Value = '6.58'
$ws.Range("E33").Value = '  +6.66%  '
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("E35").Value = '  +5.06%  '
$ws.Range("D36").Value = '26.33'
$ws.Range("E36").Value = '  +2.17%  '
$ws.Range("D37").Value = '3.619.79'
$ws.Range("E37").Value = '  +2.43%  '
$ws.Range("D38").Value = '8.34'
$ws.Range("E38").Value = '  +3.86%  '
$ws.Range("E39").Value = '  +9.31%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").Value = '0.0924'
$ws.Range("E41").Value = '  +6.75%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = '175.88'
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").Value = '31.48'
$ws.Range("E45").Value = '  +14.12%  '
$ws.Range("D46").Value = '0.912'
$ws.Range("E46").Value = '  +2.13%  '
$ws.Range("E47").Value = '  +12.10%  '
$ws.Range("D48").Value = '2.81'
$ws.Range("E48").Value = '  +7.18%  '
$ws.Range("D49").Value = '46.21'
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").Value = '7.77'
$ws.Range("E50").Value = '  +2.47%  '
$ws.Range("D51").Value = '0.266'
$ws.Range("E51").Value = '  +6.78%  '
